# Auto-generated edit script applying numeric updates to Aegis_Profits workbook sheets
# Each block targets one worksheet and updates specific cells per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Range("H116").Value = 2159.1667
$ws.Range("I116").Value = 1395
$ws.Range("J116").Value = 5980
$ws.Range("K116").Value = 1395
$ws.Range("L116").Value = 5980
$ws.Range("M116").Value = 2047
$ws.Range("N116").Value = -12864
# Row 131
$ws.Range("H131").Value = 2726.5
$ws.Range("I131").Value = 2183.25
$ws.Range("J131").Value = 4899.5
$ws.Range("K131").Value = 6549.75
$ws.Range("L131").Value = 14698.5
$ws.Range("M131").Value = -1509.75
$ws.Range("N131").Value = -24778.5
# Row 132
$ws.Range("H132").Value = 13897017
$ws.Range("I132").Value = 14714371
$ws.Range("K132").Value = 44143113
$ws.Range("M132").Value = -44140583
# Row 137
$ws.Range("H137").Value = 1338.1111
$ws.Range("I137").Value = 870.7895
$ws.Range("J137").Value = 1591.8
$ws.Range("K137").Value = 2612.3685
$ws.Range("L137").Value = 4775.4
$ws.Range("M137").Value = -62.36850000000004
$ws.Range("N137").Value = -9875.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 31003.346
$ws.Range("I32").Value = 5586.976
$ws.Range("K32").Value = 5586.976
$ws.Range("M32").Value = -5299.976
# Row 61
$ws.Range("H61").Value = 2201.0977
$ws.Range("I61").Value = 1432.2142
$ws.Range("J61").Value = 2599.7778
$ws.Range("K61").Value = 1432.2142
$ws.Range("L61").Value = 2599.7778
$ws.Range("M61").Value = -1220.2142
$ws.Range("N61").Value = -3023.7778
# Row 74
$ws.Range("H74").Value = 1770.0769
$ws.Range("I74").Value = 1829.8
$ws.Range("J74").Value = 1732.75
$ws.Range("K74").Value = 1829.8
$ws.Range("L74").Value = 1732.75
$ws.Range("M74").Value = -955.8
$ws.Range("N74").Value = -3480.75
# Row 77
$ws.Range("H77").Value = 1770.0769
$ws.Range("I77").Value = 1829.8
$ws.Range("J77").Value = 1732.75
$ws.Range("K77").Value = 9149
$ws.Range("L77").Value = 8663.75
$ws.Range("M77").Value = -4781
$ws.Range("N77").Value = -17399.75
# Row 80
$ws.Range("H80").Value = 23571
$ws.Range("J80").Value = 27485.2
$ws.Range("L80").Value = 27485.2
$ws.Range("N80").Value = -29481.2
# Row 83
$ws.Range("H83").Value = 23571
$ws.Range("J83").Value = 27485.2
$ws.Range("L83").Value = 82455.60000000001
$ws.Range("N83").Value = -92439.60000000001
# Row 110
$ws.Range("H110").Value = 26375298
$ws.Range("I110").Value = 33407516
$ws.Range("K110").Value = 33407516
$ws.Range("M110").Value = -33405471
# Row 136
$ws.Range("H136").Value = 2201.0977
$ws.Range("I136").Value = 1432.2142
$ws.Range("J136").Value = 2599.7778
$ws.Range("K136").Value = 4296.642599999999
$ws.Range("L136").Value = 7799.3334
$ws.Range("M136").Value = -1746.642599999999
$ws.Range("N136").Value = -12899.3334

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 13957.412
$ws.Range("I82").Value = 3132.182
$ws.Range("K82").Value = 3132.182
$ws.Range("M82").Value = -2749.182
# Row 85
$ws.Range("H85").Value = 13957.412
$ws.Range("I85").Value = 3132.182
$ws.Range("K85").Value = 3132.182
$ws.Range("M85").Value = -1806.182
# Row 107
$ws.Range("H107").Value = 15153586
$ws.Range("I107").Value = 20834556
$ws.Range("K107").Value = 20834556
$ws.Range("M107").Value = -20832636
# Row 134
$ws.Range("H134").Value = 1720.4736
$ws.Range("I134").Value = 1867.8928
$ws.Range("J134").Value = 1307.7
$ws.Range("K134").Value = 5603.678400000001
$ws.Range("L134").Value = 3923.1
$ws.Range("M134").Value = -3068.678400000001
$ws.Range("N134").Value = -8993.1

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 197.09525
$ws.Range("I7").Value = 61.727272
$ws.Range("J7").Value = 346
$ws.Range("K7").Value = 61.727272
$ws.Range("L7").Value = 346
$ws.Range("M7").Value = 51.272728
$ws.Range("N7").Value = -572
# Row 31
$ws.Range("H31").Value = 16391.82
$ws.Range("I31").Value = 30361.97
$ws.Range("J31").Value = 1998.3334
$ws.Range("K31").Value = 30361.97
$ws.Range("L31").Value = 1998.3334
$ws.Range("M31").Value = -30066.97
$ws.Range("N31").Value = -2588.3334
# Row 34
$ws.Range("H34").Value = 16391.82
$ws.Range("I34").Value = 30361.97
$ws.Range("J34").Value = 1998.3334
$ws.Range("K34").Value = 30361.97
$ws.Range("L34").Value = 1998.3334
$ws.Range("M34").Value = -30159.97
$ws.Range("N34").Value = -2402.3334
# Row 74
$ws.Range("H74").Value = 19269.889
$ws.Range("J74").Value = 19269.889
$ws.Range("L74").Value = 19269.889
$ws.Range("N74").Value = -21017.889
# Row 77
$ws.Range("H77").Value = 19269.889
$ws.Range("J77").Value = 19269.889
$ws.Range("L77").Value = 57809.667
$ws.Range("N77").Value = -66545.667
# Row 132
$ws.Range("H132").Value = 5413.364
$ws.Range("I132").Value = 7284.8
$ws.Range("K132").Value = 21854.4
$ws.Range("M132").Value = -19324.4

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 19.133333
$ws.Range("J2").Value = 19.76923
$ws.Range("L2").Value = 118.61538
$ws.Range("N2").Value = -344.61538
# Row 92
$ws.Range("H92").Value = 436
$ws.Range("I92").Value = 538.3333
$ws.Range("J92").Value = 282.5
$ws.Range("K92").Value = 1614.9999
$ws.Range("L92").Value = 847.5
$ws.Range("M92").Value = -366.9999
$ws.Range("N92").Value = -3343.5
# Row 131
$ws.Range("H131").Value = 3061.6667
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 3061.6667
$ws.Range("K131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("M131").Value = 9185.000100000001
$ws.Range("N131").Value = -19265.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 222.8421
$ws.Range("I2").Value = 202.83333
$ws.Range("J2").Value = 257.14285
$ws.Range("K2").Value = 202.83333
$ws.Range("L2").Value = 257.14285
$ws.Range("M2").Value = -89.83332999999999
$ws.Range("N2").Value = -483.14285
# Row 46
$ws.Range("H46").Value = 9059.6
$ws.Range("I46").Value = 6500
$ws.Range("J46").Value = 10766
$ws.Range("K46").Value = 6500
$ws.Range("L46").Value = 10766
$ws.Range("M46").Value = -6344
$ws.Range("N46").Value = -11078
# Row 57
$ws.Range("H57").Value = 9700
$ws.Range("I57").Value = 500
$ws.Range("J57").Value = 18900
$ws.Range("K57").Value = 500
$ws.Range("L57").Value = 18900
$ws.Range("M57").Value = 320
$ws.Range("N57").Value = -20540
# Row 80
$ws.Range("H80").Value = 5650
$ws.Range("I80").Value = 8800
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 8800
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -7802
$ws.Range("N80").Value = -4496
# Row 83
$ws.Range("H83").Value = 5650
$ws.Range("I83").Value = 8800
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 44000
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -39008
$ws.Range("N83").Value = -22484
# Row 132
$ws.Range("H132").Value = 3791.6667
$ws.Range("I132").Value = 3468.6667
$ws.Range("J132").Value = 4437.6665
$ws.Range("K132").Value = 10406.0001
$ws.Range("L132").Value = 13312.9995
$ws.Range("M132").Value = -7876.000100000001
$ws.Range("N132").Value = -18372.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 42541.832
$ws.Range("I16").Value = 53242.26
$ws.Range("J16").Value = 1880.2
$ws.Range("K16").Value = 53242.26
$ws.Range("L16").Value = 1880.2
$ws.Range("M16").Value = -53072.26
$ws.Range("N16").Value = -2220.2
# Row 22
$ws.Range("H22").Value = 787.1070999999999
$ws.Range("I22").Value = 1099.25
$ws.Range("K22").Value = 1099.25
$ws.Range("M22").Value = -804.25
# Row 27
$ws.Range("H27").Value = 787.1070999999999
$ws.Range("I27").Value = 1099.25
$ws.Range("K27").Value = 1099.25
$ws.Range("M27").Value = -992.25
# Row 46
$ws.Range("H46").Value = 3099.9
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 3333.3333
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 3333.3333
$ws.Range("M46").Value = -811
$ws.Range("N46").Value = -3709.3333
# Row 68
$ws.Range("H68").Value = 4646.3335
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 4646.3335
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3889.0417
$ws.Range("I132").Value = 5661.1816
$ws.Range("J132").Value = 2389.5386
$ws.Range("K132").Value = 16983.5448
$ws.Range("L132").Value = 7168.6158
$ws.Range("M132").Value = -14453.5448
$ws.Range("N132").Value = -12228.6158
